$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 16.820675
$ws.Range("H2").Value = 50.462025
$ws.Range("I2").Value = 0.8427583848046372
$ws.Range("J2").Value = 0.8427583848046373
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 35.31114333333333
$ws.Range("N2").Value = 105.93343
$ws.Range("O2").Value = 0.6187867769880316
$ws.Range("P2").Value = 0.6187867769880316
$ws.Range("Q2").Value = 593.9572658884166
$ws.Range("R2").Value = 5345.61539299575
$ws.Range("S2").Value = 0.5214877447129007
$ws.Range("T2").Value = 0.5214877447129008

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 16.820675
$ws.Range("H3").Value = 50.462025
$ws.Range("I3").Value = 0.8427583848046372
$ws.Range("J3").Value = 0.8427583848046373
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 6.551362
$ws.Range("N3").Value = 19.654086
$ws.Range("O3").Value = 0.1148050103785518
$ws.Range("P3").Value = 0.1148050103785518
$ws.Range("Q3").Value = 110.19833100935
$ws.Range("R3").Value = 991.7849790841499
$ws.Range("S3").Value = 0.09675288511410794
$ws.Range("T3").Value = 0.09675288511410796

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 16.820675
$ws.Range("H4").Value = 50.462025
$ws.Range("I4").Value = 0.8427583848046372
$ws.Range("J4").Value = 0.8427583848046373
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.951915
$ws.Range("N4").Value = 32.855745
$ws.Range("O4").Value = 0.191919590955288
$ws.Range("P4").Value = 0.191919590955288
$ws.Range("Q4").Value = 184.218602842625
$ws.Range("R4").Value = 1657.967425583625
$ws.Range("S4").Value = 0.1617418444858452
$ws.Range("T4").Value = 0.1617418444858452

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 16.820675
$ws.Range("H5").Value = 50.462025
$ws.Range("I5").Value = 0.8427583848046372
$ws.Range("J5").Value = 0.8427583848046373
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.250702333333333
$ws.Range("N5").Value = 12.752107
$ws.Range("O5").Value = 0.07448862167812857
$ws.Range("P5").Value = 0.07448862167812857
$ws.Range("Q5").Value = 71.49968247074165
$ws.Range("R5").Value = 643.4971422366749
$ws.Range("S5").Value = 0.06277591049178331
$ws.Range("T5").Value = 0.06277591049178333

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.464483
$ws.Range("H6").Value = 4.393449
$ws.Range("I6").Value = 0.07337430439942808
$ws.Range("J6").Value = 0.07337430439942808
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 35.31114333333333
$ws.Range("N6").Value = 105.93343
$ws.Range("O6").Value = 0.6187867769880316
$ws.Range("P6").Value = 0.6187867769880316
$ws.Range("Q6").Value = 51.71256912223001
$ws.Range("R6").Value = 465.41312210007
$ws.Range("S6").Value = 0.04540304933306085
$ws.Range("T6").Value = 0.04540304933306085

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.464483
$ws.Range("H7").Value = 4.393449
$ws.Range("I7").Value = 0.07337430439942808
$ws.Range("J7").Value = 0.07337430439942808
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 6.551362
$ws.Range("N7").Value = 19.654086
$ws.Range("O7").Value = 0.1148050103785518
$ws.Range("P7").Value = 0.1148050103785518
$ws.Range("Q7").Value = 9.594358275846002
$ws.Range("R7").Value = 86.34922448261401
$ws.Range("S7").Value = 0.008423737778095361
$ws.Range("T7").Value = 0.008423737778095361

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.464483
$ws.Range("H8").Value = 4.393449
$ws.Range("I8").Value = 0.07337430439942808
$ws.Range("J8").Value = 0.07337430439942808
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 10.951915
$ws.Range("N8").Value = 32.855745
$ws.Range("O8").Value = 0.191919590955288
$ws.Range("P8").Value = 0.191919590955288
$ws.Range("Q8").Value = 16.038893334945
$ws.Range("R8").Value = 144.350040014505
$ws.Range("S8").Value = 0.01408196648696702
$ws.Range("T8").Value = 0.01408196648696702

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.464483
$ws.Range("H9").Value = 4.393449
$ws.Range("I9").Value = 0.07337430439942808
$ws.Range("J9").Value = 0.07337430439942808
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.250702333333333
$ws.Range("N9").Value = 12.752107
$ws.Range("O9").Value = 0.07448862167812857
$ws.Range("P9").Value = 0.07448862167812857
$ws.Range("Q9").Value = 6.225081305227
$ws.Range("R9").Value = 56.025731747043
$ws.Range("S9").Value = 0.005465550801304842
$ws.Range("T9").Value = 0.005465550801304842

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.1122456666666666
$ws.Range("H10").Value = 0.336737
$ws.Range("I10").Value = 0.005623791954919746
$ws.Range("J10").Value = 0.005623791954919746
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 35.31114333333333
$ws.Range("N10").Value = 105.93343
$ws.Range("O10").Value = 0.6187867769880316
$ws.Range("P10").Value = 0.6187867769880316
$ws.Range("Q10").Value = 3.963522824212221
$ws.Range("R10").Value = 35.67170541791
$ws.Range("S10").Value = 0.003479928098236012
$ws.Range("T10").Value = 0.003479928098236012

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.1122456666666666
$ws.Range("H11").Value = 0.336737
$ws.Range("I11").Value = 0.005623791954919746
$ws.Range("J11").Value = 0.005623791954919746
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 6.551362
$ws.Range("N11").Value = 19.654086
$ws.Range("O11").Value = 0.1148050103785518
$ws.Range("P11").Value = 0.1148050103785518
$ws.Range("Q11").Value = 0.7353619952646665
$ws.Range("R11").Value = 6.618257957381999
$ws.Range("S11").Value = 0.0006456394937513777
$ws.Range("T11").Value = 0.0006456394937513777

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.1122456666666666
$ws.Range("H12").Value = 0.336737
$ws.Range("I12").Value = 0.005623791954919746
$ws.Range("J12").Value = 0.005623791954919746
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 10.951915
$ws.Range("N12").Value = 32.855745
$ws.Range("O12").Value = 0.191919590955288
$ws.Range("P12").Value = 0.191919590955288
$ws.Range("Q12").Value = 1.229305000451666
$ws.Range("R12").Value = 11.063745004065
$ws.Range("S12").Value = 0.001079315851605837
$ws.Range("T12").Value = 0.001079315851605837

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.1122456666666666
$ws.Range("H13").Value = 0.336737
$ws.Range("I13").Value = 0.005623791954919746
$ws.Range("J13").Value = 0.005623791954919746
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 4.250702333333333
$ws.Range("N13").Value = 12.752107
$ws.Range("O13").Value = 0.07448862167812857
$ws.Range("P13").Value = 0.07448862167812857
$ws.Range("Q13").Value = 0.4771229172065554
$ws.Range("R13").Value = 4.294106254858999
$ws.Range("S13").Value = 0.00041890851132652
$ws.Range("T13").Value = 0.00041890851132652

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 1.561668
$ws.Range("H14").Value = 4.685003999999999
$ws.Range("I14").Value = 0.07824351884101489
$ws.Range("J14").Value = 0.07824351884101489
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 35.31114333333333
$ws.Range("N14").Value = 105.93343
$ws.Range("O14").Value = 0.6187867769880316
$ws.Range("P14").Value = 0.6187867769880316
$ws.Range("Q14").Value = 55.14428258707999
$ws.Range("R14").Value = 496.2985432837199
$ws.Range("S14").Value = 0.04841605484383393
$ws.Range("T14").Value = 0.04841605484383393

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 1.561668
$ws.Range("H15").Value = 4.685003999999999
$ws.Range("I15").Value = 0.07824351884101489
$ws.Range("J15").Value = 0.07824351884101489
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 6.551362
$ws.Range("N15").Value = 19.654086
$ws.Range("O15").Value = 0.1148050103785518
$ws.Range("P15").Value = 0.1148050103785518
$ws.Range("Q15").Value = 10.231052391816
$ws.Range("R15").Value = 92.07947152634398
$ws.Range("S15").Value = 0.00898274799259713
$ws.Range("T15").Value = 0.00898274799259713

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 1.561668
$ws.Range("H16").Value = 4.685003999999999
$ws.Range("I16").Value = 0.07824351884101489
$ws.Range("J16").Value = 0.07824351884101489
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 10.951915
$ws.Range("N16").Value = 32.855745
$ws.Range("O16").Value = 0.191919590955288
$ws.Range("P16").Value = 0.191919590955288
$ws.Range("Q16").Value = 17.10325519422
$ws.Range("R16").Value = 153.92929674798
$ws.Range("S16").Value = 0.01501646413086995
$ws.Range("T16").Value = 0.01501646413086995

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 1.561668
$ws.Range("H17").Value = 4.685003999999999
$ws.Range("I17").Value = 0.07824351884101489
$ws.Range("J17").Value = 0.07824351884101489
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 4.250702333333333
$ws.Range("N17").Value = 12.752107
$ws.Range("O17").Value = 0.07448862167812857
$ws.Range("P17").Value = 0.07448862167812857
$ws.Range("Q17").Value = 6.638185811491999
$ws.Range("R17").Value = 59.74367230342799
$ws.Range("S17").Value = 0.005828251873713884
$ws.Range("T17").Value = 0.005828251873713884
